# Summer 2014 AR.xlsx
#
# Fixed bug where classes were being unaccounted for because of my silly
# way of delimitting. Fixed bug where stat would not show up because STAT
# would show up in another School and I wasn't resetting the course.
#
# Three course blocks (header row + data row + blank separator row) were
# dropped from the "COSC" section by the old parser: COSC-275 / BIGELOW B,
# COSC-301 / WILLIAMSON K and COSC-321 / CHOI K. They belong right after
# COSC-254 and before COSC-325, so we insert 9 rows there (3 blocks * 3
# rows) and fill them in; every row below shifts down by 9, which is why
# the sheet's used range grows from A1:H153 to A1:H162.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value (e.g. "19.23%") into a cell without
# letting Excel's input-parsing turn it into a percentage-formatted
# number - mark the cell as Text first, assign, then drop back to the
# workbook's default ("Normal") style so no stray formatting is left
# behind on the cell.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

# Make room: insert 9 blank rows right before the old row 50 (COSC-325),
# shifting COSC-325 and everything after it down by 9 rows.
$ws.Rows("50:58").Insert()

# --- COSC-275 / BIGELOW B ---------------------------------------------
$ws.Range("A50").Value2 = "COSC-275"
$ws.Range("B51").Value2 = "BIGELOW B"
$ws.Range("C51").Value2 = 2.692
Set-TextValue $ws.Range("D51") "19.23%"
Set-TextValue $ws.Range("E51") "50.00%"
Set-TextValue $ws.Range("F51") "19.23%"
Set-TextValue $ws.Range("G51") "3.85%"
Set-TextValue $ws.Range("H51") "7.69%"

# --- COSC-301 / WILLIAMSON K -------------------------------------------
$ws.Range("A53").Value2 = "COSC-301"
$ws.Range("B54").Value2 = "WILLIAMSON K"
$ws.Range("C54").Value2 = 3.158
Set-TextValue $ws.Range("D54") "24.56%"
Set-TextValue $ws.Range("E54") "66.67%"
Set-TextValue $ws.Range("F54") "8.77%"
Set-TextValue $ws.Range("G54") "0.00%"
Set-TextValue $ws.Range("H54") "0.00%"

# --- COSC-321 / CHOI K --------------------------------------------------
$ws.Range("A56").Value2 = "COSC-321"
$ws.Range("B57").Value2 = "CHOI K"
$ws.Range("C57").Value2 = 2.63
Set-TextValue $ws.Range("D57") "25.93%"
Set-TextValue $ws.Range("E57") "25.93%"
Set-TextValue $ws.Range("F57") "40.74%"
Set-TextValue $ws.Range("G57") "0.00%"
Set-TextValue $ws.Range("H57") "7.41%"
